$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (Förändrad / "Changed" date) holds serial date 45178 (2023-09-09)
# for every data row (2..89). Bump it to 45179 (2023-09-10) for all rows.
for ($r = 2; $r -le 89; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value = 45179
    }
}
